$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table formatting touch-ups on the confidence-scale table:
#    - normalize the table's preferred width (50%) representation
#    - mark the first row as a repeating table header row
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$tbl.PreferredWidth = 125

$headerRow = $tbl.Rows.Item(1)
$headerRow.HeadingFormat = $true

# ---------------------------------------------------------------------------
# 2) Text fix: drop the stray "Sketch the t-distribution..." sentence that
#    was left over from an earlier draft of the assignment.
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "Find the P-value and compare it to the level of significance. Sketch the t-distribution using the t-distribution applet.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Find the P-value and compare it to the level of significance.",
    2)
